$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The last paragraph of the document ("Hecho <random-forests link>.") carries
# the "_GoBack" bookmark at its very end. That bookmark needs to move to the
# new final paragraph we are about to add, so first drop it from its current
# spot (this only removes the bookmark markers, not any text).
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# Append the new content at the end of the document:
#   (existing empty paragraph stays)
#   13/08/2020
#   Hecho <link to the "edit" exercise>.
# Placeholder tokens mark where the hyperlink and the bookmark must go; they
# get resolved/removed in the following steps.
# ---------------------------------------------------------------------------
$end = $d.Content.End
$tail = $d.Range($end, $end)
$tail.InsertAfter("`r13/08/2020`rHecho XXHYPERLINKXX.ZZBOOKMARKZZ")

# ---------------------------------------------------------------------------
# Turn the XXHYPERLINKXX placeholder into a real hyperlink run, styled like
# every other link in this document ("Hipervnculo").
# ---------------------------------------------------------------------------
$linkRange = $d.Content
$linkRange.Find.Execute("XXHYPERLINKXX", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$url = "https://www.kaggle.com/felipeescaleragonz/exercise-machine-learning-competitions/edit"
$d.Hyperlinks.Add($linkRange, $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null
$newHyperlink = $d.Hyperlinks.Item($d.Hyperlinks.Count)
$newHyperlink.Range.Style = "Hipervnculo"

# ---------------------------------------------------------------------------
# Re-create the "_GoBack" bookmark right after the trailing "." of the new
# paragraph (i.e. exactly where the ZZBOOKMARKZZ placeholder sits), then wipe
# the placeholder text away again.
# ---------------------------------------------------------------------------
$markRange = $d.Content
$markRange.Find.Execute("ZZBOOKMARKZZ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$bookmarkSpot = $d.Range($markRange.Start, $markRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot) | Out-Null

$cleanupRange = $d.Content
$cleanupRange.Find.Execute("ZZBOOKMARKZZ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 2) | Out-Null
